$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old row 2 (shifts all subsequent rows up by one)
$ws.Rows.Item(2).Delete()

# Append new rows of data at the bottom (rows 21-31)
$ws.Cells.Item(21, 1).Value = 6.207096099853516
$ws.Cells.Item(21, 2).Value = -14.30157470703125
$ws.Cells.Item(21, 3).Value = 1.141430854797363
$ws.Cells.Item(22, 1).Value = 2.070674419403076
$ws.Cells.Item(22, 2).Value = -3.818307399749756
$ws.Cells.Item(22, 3).Value = 2.332929134368896
$ws.Cells.Item(23, 1).Value = 12.43622970581055
$ws.Cells.Item(23, 2).Value = 1.804691195487976
$ws.Cells.Item(23, 3).Value = -0.571514368057251
$ws.Cells.Item(24, 1).Value = -4.415188312530518
$ws.Cells.Item(24, 2).Value = 2.631336450576782
$ws.Cells.Item(24, 3).Value = -0.9896306991577148
$ws.Cells.Item(25, 1).Value = -4.100935459136963
$ws.Cells.Item(25, 2).Value = 1.459545493125916
$ws.Cells.Item(25, 3).Value = 7.553257465362549
$ws.Cells.Item(26, 1).Value = -3.94806981086731
$ws.Cells.Item(26, 2).Value = 8.382166862487793
$ws.Cells.Item(26, 3).Value = -2.697782278060913
$ws.Cells.Item(27, 1).Value = 4.418517112731934
$ws.Cells.Item(27, 2).Value = 15.49866580963135
$ws.Cells.Item(27, 3).Value = -7.229417324066162
$ws.Cells.Item(28, 1).Value = 1.811815142631531
$ws.Cells.Item(28, 2).Value = -5.760817050933838
$ws.Cells.Item(28, 3).Value = -5.106345176696777
$ws.Cells.Item(29, 1).Value = 6.468618869781494
$ws.Cells.Item(29, 2).Value = 1.199088335037231
$ws.Cells.Item(29, 3).Value = -6.144978046417236
$ws.Cells.Item(30, 1).Value = -2.402371168136597
$ws.Cells.Item(30, 2).Value = 2.291517019271851
$ws.Cells.Item(30, 3).Value = 1.192030906677246
$ws.Cells.Item(31, 1).Value = -12.0989408493042
$ws.Cells.Item(31, 2).Value = -13.58198833465576
$ws.Cells.Item(31, 3).Value = -2.806971788406372
